# Rename the embedded logo pictures (InlineShape "Name" -> wp:docPr/@name):
#   - the two Pearson Edexcel logos that live in the document's footers
#     go from "image2.png" to "image1.png"
#   - the BTEC logo that lives in the document's (first-page) header
#     goes from "image1.jpg" to "image2.jpg"
#
# Renaming is done by round-tripping InlineShape -> Shape -> InlineShape;
# setting .Name directly on some InlineShapes in header/footer stories can
# otherwise fail to resolve, so the Shape detour is used everywhere for
# consistency and reliability.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlineShape($inlineShape, $newName) {
    $asShape = $inlineShape.ConvertToShape()
    $asShape.Name = $newName
    $asShape.ConvertToInlineShape() | Out-Null
}

# --- Footers: Pearson Edexcel logo, in both the default and first-page
#     footers of section 1 ---
for ($f = 1; $f -le $sec.Footers.Count; $f++) {
    $ftr = $sec.Footers.Item($f)
    if ($ftr.Exists) {
        $shapes = $ftr.Range.InlineShapes
        for ($i = 1; $i -le $shapes.Count; $i++) {
            $shp = $shapes.Item($i)
            if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                Rename-InlineShape $shp "image1.png"
            }
        }
    }
}

# --- Headers: BTEC logo, wherever it appears (the first-page header) ---
for ($h = 1; $h -le $sec.Headers.Count; $h++) {
    $hdr = $sec.Headers.Item($h)
    if ($hdr.Exists) {
        $shapes = $hdr.Range.InlineShapes
        for ($i = 1; $i -le $shapes.Count; $i++) {
            $shp = $shapes.Item($i)
            if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                Rename-InlineShape $shp "image2.jpg"
            }
        }
    }
}

Write-Output "Done renaming logo InlineShapes."
